# Adds a new "2022-Q3" quarter: a new sheet with fund-holdings data, and a
# new summary row at the top of "总计" (with every later summary row shifted
# down by one).

$wb = $excel.ActiveWorkbook

function NeedsPrefix($s) {
    # Values that look like a bare number must get a leading apostrophe so
    # Excel stores them as text (matching the source file's inlineStr cells)
    # instead of silently coercing them into numbers (which would also
    # destroy leading zeros in fund codes like "013220").
    if ($s -match '^-?\d+(\.\d+)?$') {
        return $true
    }
    return $false
}

function SetTextValue($cell, $s) {
    if (NeedsPrefix $s) {
        $cell.Value = "'" + $s
    } else {
        $cell.Value = $s
    }
}

# ---------------------------------------------------------------------
# Step 1: new "2022-Q3" worksheet (fund holdings), inserted right after
# "总计". Copying "2022-Q2" gives us the exact same header row/column
# styling, then we overwrite the data and drop the rows we don't need.
# ---------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("2022-Q2")
$afterSheet = $wb.Worksheets.Item("总计")
$srcSheet.Copy($null, $afterSheet)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q3"

$q3Data = @(
    @(0, "166019", "中欧价值智选回报混合A", "71.01", "94.44", "7.87", "5.5885", 2),
    @(1, "013220", "中欧新兴价值一年持有混合A", "34.31", "94.30", "8.03", "2.7551", 3),
    @(2, "004235", "中欧价值智选回报混合C", "26.97", "94.44", "7.87", "2.1225", 2),
    @(3, "013221", "中欧新兴价值一年持有混合C", "15.22", "94.30", "8.03", "1.2222", 3),
    @(4, "004848", "中欧睿泓定期开放灵活配置混合", "17.10", "59.07", "6.07", "1.0380", 2),
    @(5, "014404", "中欧多元价值三年持有混合A", "12.96", "94.15", "7.72", "1.0005", 3),
    @(6, "001887", "中欧价值智选回报混合E", "12.27", "94.44", "7.87", "0.9656", 2),
    @(7, "012568", "天弘高端制造混合A", "5.98", "92.42", "5.74", "0.3433", 4),
    @(8, "000418", "景顺长城成长之星", "6.70", "81.16", "4.38", "0.2935", 3),
    @(9, "014405", "中欧多元价值三年持有混合C", "1.23", "94.15", "7.72", "0.0950", 3),
    @(10, "168501", "北信瑞丰产业升级多策略混合", "1.61", "93.34", "5.26", "0.0847", 1),
    @(11, "009598", "景顺长城科技创新三年定期开放灵活配置混合", "2.73", "91.12", "2.77", "0.0756", 10),
    @(12, "012569", "天弘高端制造混合C", "0.97", "92.42", "5.74", "0.0557", 4),
    @(13, "005041", "人保研究精选混合A", "1.23", "79.96", "2.43", "0.0299", 8),
    @(14, "006973", "太平睿盈混合A", "3.84", "28.79", "0.75", "0.0288", 7),
    @(15, "002123", "北信瑞丰外延增长主题灵活配置混合", "0.15", "88.55", "9.04", "0.0136", 1),
    @(16, "007669", "太平睿盈混合C", "1.04", "28.79", "0.75", "0.0078", 7),
    @(17, "005042", "人保研究精选混合C", "0.08", "79.96", "2.43", "0.0019", 8)
)

$r = 2
foreach ($item in $q3Data) {
    $newSheet.Cells.Item($r, 1).Value = $item[0]
    SetTextValue $newSheet.Cells.Item($r, 2) $item[1]
    SetTextValue $newSheet.Cells.Item($r, 3) $item[2]
    SetTextValue $newSheet.Cells.Item($r, 4) $item[3]
    SetTextValue $newSheet.Cells.Item($r, 5) $item[4]
    SetTextValue $newSheet.Cells.Item($r, 6) $item[5]
    SetTextValue $newSheet.Cells.Item($r, 7) $item[6]
    $newSheet.Cells.Item($r, 8).Value = $item[7]
    $r = $r + 1
}

# 2022-Q2 had 29 data rows (30 incl. header); 2022-Q3 only needs 18 (19
# incl. header), so drop the leftover rows copied from the source sheet.
$newSheet.Range("A20:H30").EntireRow.Delete()

# ---------------------------------------------------------------------
# Step 2: "总计" summary sheet - add the 2022-Q3 row at the top and shift
# every other quarter down by one row.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$sheet1Data = @(
    @(0, "2022-Q3", 18, 15.72),
    @(1, "2022-Q2", 29, 22.08),
    @(2, "2022-Q1", 35, 35.73),
    @(3, "2021-Q4", 60, 56),
    @(4, "2021-Q3", 54, 40.95),
    @(5, "2021-Q2", 67, 32.74),
    @(6, "2021-Q1", 69, 31.21),
    @(7, "2020-Q4", 95, 39.32)
)

$row = 2
foreach ($item in $sheet1Data) {
    $totalSheet.Cells.Item($row, 1).Value = $item[0]
    $totalSheet.Cells.Item($row, 2).Value = $item[1]
    $totalSheet.Cells.Item($row, 3).Value = $item[2]
    $totalSheet.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

# Row 9 (2020-Q4) is brand new - give A9 the same style as the other index
# cells above it (border + center alignment), since a freshly written cell
# otherwise gets no style at all.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)
$totalSheet.Range("A9").Value = 7

Write-Host "edit complete"
